$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (pushes the existing rows 9-10 down to 10-11,
# carrying their values/formatting with them — this reproduces the
# "old row 9 -> new row 10" and "old row 10 -> new row 11" shift in the diff).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with this week's record.
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C9").Value = "Los Lagos"

# Match the date-time number format already used by the other rows in
# column D before assigning the serial date value, so the cell reuses the
# existing style instead of Excel minting a brand-new one.
$ws.Range("D9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D9").Value = 44516

$ws.Range("E9").Value = 10
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107002
$ws.Range("J9").Value = "Chirimoya"
$ws.Range("K9").Value = "Cultivar IV Región"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 19000
$ws.Range("P9").Value = 18500
$ws.Range("Q9").Value = "$/bandeja 8 kilos"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 2312
$ws.Range("T9").Value = 8
